$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 2022.0358
$ws.Range("J40").Value = 2141.4167
$ws.Range("L40").Value = 2141.4167
$ws.Range("N40").Value = -2491.4167
# Row 64
$ws.Range("H64").Value = 102867.9
$ws.Range("I64").Value = 202396
$ws.Range("J64").Value = 3339.8
$ws.Range("K64").Value = 202396
$ws.Range("L64").Value = 3339.8
$ws.Range("M64").Value = -202148
$ws.Range("N64").Value = -3835.8
# Row 67
$ws.Range("H67").Value = 102867.9
$ws.Range("I67").Value = 202396
$ws.Range("J67").Value = 3339.8
$ws.Range("K67").Value = 202396
$ws.Range("L67").Value = 3339.8
$ws.Range("M67").Value = -201538
$ws.Range("N67").Value = -5055.8
# Row 74
$ws.Range("H74").Value = 3757.5
$ws.Range("I74").Value = 4333.3335
$ws.Range("J74").Value = 3412
$ws.Range("K74").Value = 4333.3335
$ws.Range("L74").Value = 3412
$ws.Range("M74").Value = -3397.3335
$ws.Range("N74").Value = -5284
# Row 77
$ws.Range("H77").Value = 3757.5
$ws.Range("I77").Value = 4333.3335
$ws.Range("J77").Value = 3412
$ws.Range("K77").Value = 21666.6675
$ws.Range("L77").Value = 17060
$ws.Range("M77").Value = -16986.6675
$ws.Range("N77").Value = -26420
# Row 100
$ws.Range("H100").Value = 1629
$ws.Range("I100").Value = 1782
$ws.Range("J100").Value = 1519.7142
$ws.Range("K100").Value = 1782
$ws.Range("L100").Value = 1519.7142
$ws.Range("M100").Value = -1241
$ws.Range("N100").Value = -2601.7142
# Row 127
$ws.Range("H127").Value = 2169.074
$ws.Range("J127").Value = 2383.3958
$ws.Range("L127").Value = 7150.187399999999
$ws.Range("N127").Value = -17070.1874
# Row 129
$ws.Range("H129").Value = 875392.5
$ws.Range("J129").Value = 993234.3
$ws.Range("L129").Value = 2979702.9
$ws.Range("N129").Value = -2989702.9
# Row 132
$ws.Range("H132").Value = 11370855
$ws.Range("I132").Value = 11911796
$ws.Range("K132").Value = 35735388
$ws.Range("M132").Value = -35732858
# Row 138
$ws.Range("H138").Value = 2354.279
$ws.Range("I138").Value = 2232.2666
$ws.Range("J138").Value = 2380.0564
$ws.Range("K138").Value = 6696.7998
$ws.Range("L138").Value = 7140.1692
$ws.Range("M138").Value = -1556.7998
$ws.Range("N138").Value = -17420.1692

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 33373.98
$ws.Range("I32").Value = 5945.6665
$ws.Range("K32").Value = 5945.6665
$ws.Range("M32").Value = -5658.6665
# Row 74
$ws.Range("H74").Value = 1835
$ws.Range("I74").Value = 905.5161000000001
$ws.Range("K74").Value = 905.5161000000001
$ws.Range("M74").Value = -31.51610000000005
# Row 77
$ws.Range("H77").Value = 1835
$ws.Range("I77").Value = 905.5161000000001
$ws.Range("K77").Value = 4527.5805
$ws.Range("M77").Value = -159.5805
# Row 102
$ws.Range("H102").Value = 35440.2
$ws.Range("I102").Value = 64623
$ws.Range("J102").Value = 2088.4285
$ws.Range("K102").Value = 64623
$ws.Range("L102").Value = 2088.4285
$ws.Range("M102").Value = -63001
$ws.Range("N102").Value = -5332.4285
# Row 132
$ws.Range("H132").Value = 3344.84
$ws.Range("I132").Value = 3025.625
$ws.Range("J132").Value = 3912.3333
$ws.Range("K132").Value = 9076.875
$ws.Range("L132").Value = 11736.9999
$ws.Range("M132").Value = -6546.875
$ws.Range("N132").Value = -16796.9999

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 680.25
$ws.Range("I94").Value = 581.4
$ws.Range("J94").Value = 725.1818
$ws.Range("K94").Value = 581.4
$ws.Range("L94").Value = 725.1818
$ws.Range("M94").Value = -130.4
$ws.Range("N94").Value = -1627.1818
# Row 103
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 31442.715
$ws.Range("J31").Value = 3607.3635
$ws.Range("L31").Value = 3607.3635
$ws.Range("N31").Value = -4197.363499999999
# Row 34
$ws.Range("H34").Value = 31442.715
$ws.Range("J34").Value = 3607.3635
$ws.Range("L34").Value = 3607.3635
$ws.Range("N34").Value = -4011.3635
# Row 58
$ws.Range("H58").Value = 8448.056
$ws.Range("I58").Value = 1074.258
$ws.Range("J58").Value = 54165.6
$ws.Range("K58").Value = 1074.258
$ws.Range("L58").Value = 54165.6
$ws.Range("M58").Value = -871.258
$ws.Range("N58").Value = -54571.6
# Row 103
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").ClearContents()
# Row 136
$ws.Range("H136").Value = 8448.056
$ws.Range("I136").Value = 1074.258
$ws.Range("J136").Value = 54165.6
$ws.Range("K136").Value = 3222.774
$ws.Range("L136").Value = 162496.8
$ws.Range("M136").Value = -672.7740000000003
$ws.Range("N136").Value = -167596.8

$ws = $wb.Worksheets.Item("CUL")
# Row 129
$ws.Range("H129").Value = 245844.28
$ws.Range("I129").Value = 11757.2
$ws.Range("J129").Value = 318996.5
$ws.Range("K129").Value = 35271.60000000001
$ws.Range("L129").Value = 956989.5
$ws.Range("M129").Value = -30271.60000000001
$ws.Range("N129").Value = -966989.5
# Row 131
$ws.Range("H131").Value = 744.42
$ws.Range("I131").Value = 453.4375
$ws.Range("J131").Value = 799.8452
$ws.Range("K131").Value = 1360.3125
$ws.Range("L131").Value = 2399.5356
$ws.Range("M131").Value = 3679.6875
$ws.Range("N131").Value = -12479.5356
# Row 132
$ws.Range("H132").Value = 4915.6924
$ws.Range("I132").Value = 4499.8335
$ws.Range("J132").Value = 5272.143
$ws.Range("K132").Value = 40498.5015
$ws.Range("L132").Value = 47449.287
$ws.Range("M132").Value = -37968.5015
$ws.Range("N132").Value = -52509.287

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 66427.55
$ws.Range("I70").Value = 116176
$ws.Range("J70").Value = 6729.4
$ws.Range("K70").Value = 116176
$ws.Range("L70").Value = 6729.4
$ws.Range("M70").Value = -115906
$ws.Range("N70").Value = -7269.4
# Row 73
$ws.Range("H73").Value = 66427.55
$ws.Range("I73").Value = 116176
$ws.Range("J73").Value = 6729.4
$ws.Range("K73").Value = 116176
$ws.Range("L73").Value = 6729.4
$ws.Range("M73").Value = -115240
$ws.Range("N73").Value = -8601.4
# Row 109
$ws.Range("H109").Value = 9285
$ws.Range("J109").Value = 9285
$ws.Range("L109").Value = 9285
$ws.Range("N109").Value = -11365

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 844280
$ws.Range("I46").Value = 233.33333
$ws.Range("J46").Value = 1125628.9
$ws.Range("K46").Value = 233.33333
$ws.Range("L46").Value = 1125628.9
$ws.Range("M46").Value = -45.33332999999999
$ws.Range("N46").Value = -1126004.9
# Row 68
$ws.Range("H68").Value = 4952.1665
$ws.Range("I68").Value = 2400
$ws.Range("J68").Value = 5462.6
$ws.Range("K68").Value = 2400
$ws.Range("L68").Value = 5462.6
$ws.Range("M68").Value = -1651
$ws.Range("N68").Value = -6960.6
# Row 69
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
# Row 71
$ws.Range("H71").Value = 4952.1665
$ws.Range("I71").Value = 2400
$ws.Range("J71").Value = 5462.6
$ws.Range("K71").Value = 12000
$ws.Range("L71").Value = 27313
$ws.Range("M71").Value = -8256
$ws.Range("N71").Value = -34801
# Row 72
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
# Row 93
$ws.Range("H93").Value = 1746
$ws.Range("I93").Value = 2694
$ws.Range("K93").Value = 2694
$ws.Range("M93").Value = -1446

$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 167083.33
$ws.Range("I100").Value = 250300
$ws.Range("J100").Value = 650
$ws.Range("K100").Value = 500600
$ws.Range("L100").Value = 1300
$ws.Range("M100").Value = -500059
$ws.Range("N100").Value = -2382
